$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text: volume/number and week-covering dates (shared strings in A8/C9) ---
$ws.Range("A8").Value = "Volume 30   Number  46"
$ws.Range("C9").Value = "Report Covering the Week  11/13/2023  Through  11/19/2023"

# --- Weekly crime-stat numbers refreshed for the new reporting week ---
# (style/number-format of these cells is unaffected by a plain .Value assignment)
$ws.Range("F15").Value = 3
$ws.Range("I15").Value = 16
$ws.Range("K15").Value = 33.333333333333
$ws.Range("L15").Value = 23.076923076923
$ws.Range("M15").Value = -11.111111111111
$ws.Range("N15").Value = -20
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = -40
$ws.Range("G16").Value = 23
$ws.Range("H16").Value = -43.478260869565
$ws.Range("I16").Value = 180
$ws.Range("J16").Value = 192
$ws.Range("K16").Value = -6.25
$ws.Range("L16").Value = 18.421052631578
$ws.Range("M16").Value = -15.492957746478
$ws.Range("N16").Value = -56.521739130434
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 7
$ws.Range("E17").Value = -14.285714285714
$ws.Range("F17").Value = 16
$ws.Range("G17").Value = 27
$ws.Range("H17").Value = -40.74074074074
$ws.Range("I17").Value = 261
$ws.Range("J17").Value = 238
$ws.Range("K17").Value = 9.663865546218
$ws.Range("L17").Value = 27.941176470588
$ws.Range("M17").Value = 59.146341463414
$ws.Range("N17").Value = 27.317073170731
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = -33.333333333333
$ws.Range("F18").Value = 9
$ws.Range("G18").Value = 11
$ws.Range("H18").Value = -18.181818181818
$ws.Range("I18").Value = 123
$ws.Range("J18").Value = 93
$ws.Range("K18").Value = 32.258064516129
$ws.Range("L18").Value = 24.242424242424
$ws.Range("M18").Value = -55.272727272727
$ws.Range("N18").Value = -83.265306122449
$ws.Range("C19").Value = 16
$ws.Range("D19").Value = 13
$ws.Range("E19").Value = 23.076923076923
$ws.Range("F19").Value = 67
$ws.Range("H19").Value = 24.074074074074
$ws.Range("I19").Value = 596
$ws.Range("J19").Value = 511
$ws.Range("K19").Value = 16.634050880626
$ws.Range("L19").Value = 39.252336448598
$ws.Range("M19").Value = 37.962962962963
$ws.Range("N19").Value = 52.820512820512
$ws.Range("C20").Value = 6
$ws.Range("D20").Value = 5
$ws.Range("E20").Value = 20
$ws.Range("F20").Value = 28
$ws.Range("G20").Value = 23
$ws.Range("H20").Value = 21.739130434782
$ws.Range("I20").Value = 461
$ws.Range("J20").Value = 269
$ws.Range("K20").Value = 71.375464684014
$ws.Range("L20").Value = 66.425992779783
$ws.Range("M20").Value = 146.524064171123
$ws.Range("N20").Value = -75.78781512605
$ws.Range("D21").Value = 33
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 136
$ws.Range("G21").Value = 138
$ws.Range("H21").Value = -1.449275362318
$ws.Range("I21").Value = 1641
$ws.Range("J21").Value = 1319
$ws.Range("K21").Value = 24.412433661865
$ws.Range("L21").Value = 39.303904923599
$ws.Range("M21").Value = 27.307990690457
$ws.Range("N21").Value = -55.419722901385
$ws.Range("F23").Value = 2
$ws.Range("H23").Value = -33.333333333333
$ws.Range("L23").Value = 25.581395348837
$ws.Range("M23").Value = 10.204081632653
$ws.Range("C24").Value = 19
$ws.Range("D24").Value = 23
$ws.Range("E24").Value = -17.391304347826
$ws.Range("F24").Value = 83
$ws.Range("G24").Value = 94
$ws.Range("H24").Value = -11.702127659574
$ws.Range("I24").Value = 1229
$ws.Range("J24").Value = 1029
$ws.Range("K24").Value = 19.436345966958
$ws.Range("L24").Value = 39.817974971558
$ws.Range("M24").Value = -6.397562833206
$ws.Range("C25").Value = 8
$ws.Range("D25").Value = 6
$ws.Range("E25").Value = 33.333333333333
$ws.Range("F25").Value = 37
$ws.Range("G25").Value = 41
$ws.Range("H25").Value = -9.756097560975
$ws.Range("I25").Value = 454
$ws.Range("J25").Value = 408
$ws.Range("K25").Value = 11.274509803921
$ws.Range("L25").Value = 17.312661498708
$ws.Range("M25").Value = 14.936708860759
$ws.Range("F26").Value = 3
$ws.Range("I26").Value = 29
$ws.Range("K26").Value = 11.538461538461
$ws.Range("L26").Value = 70.588235294117
$ws.Range("F27").Value = 3
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 46
$ws.Range("K27").Value = -15.217391304347
$ws.Range("L28").Value = -35.714285714285
$ws.Range("L29").Value = -36.363636363636

# --- Cells that become the "no data" placeholder (shared strings "0" / "***.*") ---
# These must render with the same General, right-aligned style (s=14 in the original
# file) as every other placeholder cell on the sheet. Simply assigning .Value with a
# numeric-looking string like "0" gets auto-coerced back to a number, so: force Text
# entry via NumberFormat "@", assign the literal value, then paste just the formats
# from an existing placeholder cell (G15) on top so the style matches exactly and no
# stray custom number format is left behind.
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "0"
$ws.Range("G15").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0"
$ws.Range("G15").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "***.*"
$ws.Range("G15").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("C23").NumberFormat = "@"
$ws.Range("C23").Value = "0"
$ws.Range("G15").Copy()
$ws.Range("C23").PasteSpecial(-4122)
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0"
$ws.Range("G15").Copy()
$ws.Range("D23").PasteSpecial(-4122)
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "***.*"
$ws.Range("G15").Copy()
$ws.Range("E23").PasteSpecial(-4122)
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0"
$ws.Range("G15").Copy()
$ws.Range("D30").PasteSpecial(-4122)
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "***.*"
$ws.Range("G15").Copy()
$ws.Range("E30").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Cells that gain real counts this week (were "no data" placeholders) ---
# Give them the same number formats used by sibling cells in their column
# (#,##0 for counts, the signed one-decimal format for % change) before writing
# the value, so the resulting style matches the rest of the column exactly.
$ws.Range("D27").NumberFormat = "#,##0"
$ws.Range("D27").Value = 1
$ws.Range("E27").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("E27").Value = -100
